# Auto-generated Excel COM-interop script to apply the Maduin_Profits update.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) for specific
# leve rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets, matching
# the refreshed Universalis market-board pull from the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1636.75
$ws.Range("I2").Value = 764.2
$ws.Range("K2").Value = 764.2
$ws.Range("M2").Value = -651.2
$ws.Range("H8").Value = 76.59999999999999
$ws.Range("J8").Value = 51.333332
$ws.Range("L8").Value = 153.999996
$ws.Range("N8").Value = -431.999996
$ws.Range("H19").Value = 1744.1111
$ws.Range("J19").Value = 2133
$ws.Range("L19").Value = 2133
$ws.Range("N19").Value = -2483
$ws.Range("H28").Value = 1519.8572
$ws.Range("I28").Value = 999.25
$ws.Range("K28").Value = 999.25
$ws.Range("M28").Value = -514.25
$ws.Range("H32").Value = 1749.125
$ws.Range("I32").Value = 1624.75
$ws.Range("J32").Value = 1873.5
$ws.Range("K32").Value = 1624.75
$ws.Range("L32").Value = 1873.5
$ws.Range("M32").Value = -1298.75
$ws.Range("N32").Value = -2525.5
$ws.Range("I40").Value = 1350
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1350
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -1175
$ws.Range("H43").Value = 8702.267
$ws.Range("I43").Value = 7887.1113
$ws.Range("J43").Value = 9925
$ws.Range("K43").Value = 7887.1113
$ws.Range("L43").Value = 9925
$ws.Range("M43").Value = -7818.1113
$ws.Range("N43").Value = -10063
$ws.Range("H62").Value = 3880.1667
$ws.Range("I62").Value = 2997.6667
$ws.Range("J62").Value = 4762.6665
$ws.Range("K62").Value = 2997.6667
$ws.Range("L62").Value = 4762.6665
$ws.Range("M62").Value = -2373.6667
$ws.Range("N62").Value = -6010.6665
$ws.Range("H64").Value = 7664.75
$ws.Range("I64").Value = 4338.6
$ws.Range("J64").Value = 13208.333
$ws.Range("K64").Value = 4338.6
$ws.Range("L64").Value = 13208.333
$ws.Range("M64").Value = -4090.6
$ws.Range("N64").Value = -13704.333
$ws.Range("H65").Value = 3880.1667
$ws.Range("I65").Value = 2997.6667
$ws.Range("J65").Value = 4762.6665
$ws.Range("K65").Value = 14988.3335
$ws.Range("L65").Value = 23813.3325
$ws.Range("M65").Value = -11868.3335
$ws.Range("N65").Value = -30053.3325
$ws.Range("H67").Value = 7664.75
$ws.Range("I67").Value = 4338.6
$ws.Range("J67").Value = 13208.333
$ws.Range("K67").Value = 4338.6
$ws.Range("L67").Value = 13208.333
$ws.Range("M67").Value = -3480.6
$ws.Range("N67").Value = -14924.333
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").ClearContents()
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").ClearContents()
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = 0
$ws.Range("H76").Value = 3600
$ws.Range("I76").Value = 3600
$ws.Range("K76").Value = 3600
$ws.Range("M76").Value = -3285
$ws.Range("H79").Value = 3600
$ws.Range("I79").Value = 3600
$ws.Range("K79").Value = 3600
$ws.Range("M79").Value = -2508
$ws.Range("H111").Value = 1941.5555
$ws.Range("I111").Value = 2750
$ws.Range("J111").Value = 1710.5714
$ws.Range("K111").Value = 8250
$ws.Range("L111").Value = 5131.7142
$ws.Range("M111").Value = -5183
$ws.Range("N111").Value = -11265.7142
$ws.Range("H118").Value = 558.4
$ws.Range("J118").Value = 1499
$ws.Range("L118").Value = 4497
$ws.Range("N118").Value = -7811
$ws.Range("H138").Value = 1979.5
$ws.Range("I138").Value = 1381.9412
$ws.Range("K138").Value = 4145.8236
$ws.Range("M138").Value = 994.1764000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4342.45
$ws.Range("I32").Value = 3267.4375
$ws.Range("J32").Value = 8642.5
$ws.Range("K32").Value = 3267.4375
$ws.Range("L32").Value = 8642.5
$ws.Range("M32").Value = -2980.4375
$ws.Range("N32").Value = -9216.5
$ws.Range("H76").Value = 50000
$ws.Range("J76").Value = 50000
$ws.Range("L76").Value = 50000
$ws.Range("N76").Value = -50676
$ws.Range("H79").Value = 50000
$ws.Range("J79").Value = 50000
$ws.Range("L79").Value = 50000
$ws.Range("N79").Value = -52340
$ws.Range("H97").Value = 1703.2
$ws.Range("I97").Value = 1170.5
$ws.Range("J97").Value = 2502.25
$ws.Range("K97").Value = 1170.5
$ws.Range("L97").Value = 2502.25
$ws.Range("M97").Value = -674.5
$ws.Range("N97").Value = -3494.25
$ws.Range("H119").Value = 49250
$ws.Range("J119").Value = 49250
$ws.Range("L119").Value = 49250
$ws.Range("N119").Value = -58926

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3100.2222
$ws.Range("I105").Value = 2601.2144
$ws.Range("K105").Value = 2601.2144
$ws.Range("M105").Value = -854.2143999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 29999
$ws.Range("J4").Value = 29999
$ws.Range("L4").Value = 29999
$ws.Range("N4").Value = -30223

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 4249.5
$ws.Range("I22").Value = 1997.5
$ws.Range("J22").Value = 6501.5
$ws.Range("K22").Value = 5992.5
$ws.Range("L22").Value = 19504.5
$ws.Range("M22").Value = -5823.5
$ws.Range("N22").Value = -19842.5
$ws.Range("H27").Value = 4249.5
$ws.Range("I27").Value = 1997.5
$ws.Range("J27").Value = 6501.5
$ws.Range("K27").Value = 5992.5
$ws.Range("L27").Value = 19504.5
$ws.Range("M27").Value = -5890.5
$ws.Range("N27").Value = -19708.5
$ws.Range("H80").Value = 7331
$ws.Range("J80").Value = 7663.6665
$ws.Range("L80").Value = 22990.9995
$ws.Range("N80").Value = -24862.9995
$ws.Range("H83").Value = 7331
$ws.Range("J83").Value = 7663.6665
$ws.Range("L83").Value = 68972.9985
$ws.Range("N83").Value = -78332.9985
$ws.Range("H86").Value = 725.5
$ws.Range("I86").Value = 725.5
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2176.5
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -990.5
$ws.Range("H89").Value = 725.5
$ws.Range("I89").Value = 725.5
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 6529.5
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -601.5
$ws.Range("H112").Value = 39995.168
$ws.Range("I112").Value = 19987.5
$ws.Range("K112").Value = 59962.5
$ws.Range("M112").Value = -58854.5
$ws.Range("H116").Value = 850
$ws.Range("J116").Value = 850
$ws.Range("L116").Value = 2550
$ws.Range("N116").Value = -9434
$ws.Range("H118").Value = 229
$ws.Range("I118").Value = 229
$ws.Range("K118").Value = 687
$ws.Range("M118").Value = 556
$ws.Range("H122").Value = 233.85715
$ws.Range("J122").Value = 525.5
$ws.Range("L122").Value = 4729.5
$ws.Range("N122").Value = -9629.5
$ws.Range("H128").Value = 277692.34
$ws.Range("I128").Value = 277692.34
$ws.Range("K128").Value = 833077.02
$ws.Range("M128").Value = -828097.02
$ws.Range("H129").Value = 1603.75
$ws.Range("I129").Value = 971.6667
$ws.Range("J129").Value = 3500
$ws.Range("K129").Value = 2915.0001
$ws.Range("L129").Value = 10500
$ws.Range("M129").Value = 2084.9999
$ws.Range("N129").Value = -20500
$ws.Range("H131").Value = 971.82355
$ws.Range("I131").Value = 928.2
$ws.Range("J131").Value = 990
$ws.Range("K131").Value = 2784.6
$ws.Range("L131").Value = 2970
$ws.Range("M131").Value = 2255.4
$ws.Range("N131").Value = -13050

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1895.5745
$ws.Range("I132").Value = 1714.5111
$ws.Range("J132").Value = 5969.5
$ws.Range("K132").Value = 5143.5333
$ws.Range("L132").Value = 17908.5
$ws.Range("M132").Value = -2613.5333
$ws.Range("N132").Value = -22968.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2056.2727
$ws.Range("I16").Value = 1978
$ws.Range("J16").Value = 2265
$ws.Range("K16").Value = 1978
$ws.Range("L16").Value = 2265
$ws.Range("M16").Value = -1808
$ws.Range("N16").Value = -2605
$ws.Range("H22").Value = 1500
$ws.Range("I22").Value = 1500
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1205
$ws.Range("H27").Value = 1500
$ws.Range("I27").Value = 1500
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1500
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -1393

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 34151.332
$ws.Range("J63").Value = 34151.332
$ws.Range("L63").Value = 34151.332
$ws.Range("N63").Value = -35399.332
$ws.Range("H66").Value = 34151.332
$ws.Range("J66").Value = 34151.332
$ws.Range("L66").Value = 102453.996
$ws.Range("N66").Value = -108693.996
$ws.Range("H103").Value = 11900.5
$ws.Range("J103").Value = 11900.5
$ws.Range("L103").Value = 11900.5
$ws.Range("N103").Value = -14244.5
